$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename the "on" mode value to "trim" (all data rows BX3:BX7) -------------
$ws.Range("BX3:BX7").Value = "trim"

# 2) Update the barrnap_mode column description (BX1) with the new wording ---
$desc = "Should barrnap be run on the sample to detect the presence of rRNA genes? This column is optional. This value can be:`n" +
        "      - `"off`": Will skip gene detection (default). `n" +
        "      - `"filter`": To remove reads that didn't have hits for both 16S and 23S.`n" +
        "      - `"concat`": To remove the ITS region between 16S and 23S.`n" +
        "      - `"trim`": To remove the 23S and keep the 16S portion only of every read.`n"
$ws.Range("BX1").Value = $desc

# 3) Move the active selection from BX1 to BX9 --------------------------------
$ws.Range("BX9").Select()

# 4) Widen column BX (76) so the long description text is readable ----------
#    Column BW (75) keeps its original width; only BX (76) grows.
$ws.Columns.Item(76).ColumnWidth = 38
